$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers I1 "I0" and J1 "IF" - copy formatting (style index 1, bold/centered/border)
# from the existing header cell H1, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2..37: I = 1 (constant), J = same value as H (existing column)
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
